# Apply crypto price/volume updates
# Commit: "Updated cryptos list on Fri Jan 26 03:45:40 UTC 2024 with GitHub Actions"
#
# For Price (column D) cells whose new value looks like a plain decimal number,
# force the cell's number format to Text first so Excel keeps it as a literal
# string (matching the source data, which mixes thousand-separated values like
# "40.137.63" with plain decimals like "293.91") instead of silently coercing
# it into a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.137.63'
$ws.Range('D3').Value = '2.222.97'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.91'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.76'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.74'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '50.44'
$ws.Range('E11').Value = '  +5.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0781'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.113'
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '2.582.75'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.86'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').Value = '2.210.35'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').Value = '40.072.37'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.28'
$ws.Range('E21').Value = '  -5.45%  '
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.73'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.63'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.83'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.20'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.90'
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.82'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.97'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.02'
$ws.Range('E35').Value = '  +6.33%  '
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.32'
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.75'
$ws.Range('E39').Value = '  +2.59%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0996'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').Value = '2.084.03'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.06'
$ws.Range('E44').Value = '  +7.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.11'
$ws.Range('E45').Value = '  +2.40%  '
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.73'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  -10.69%  '
$ws.Range('D49').Value = '2.449.87'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('E50').Value = '  +2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.12'
$ws.Range('E51').Value = '  +3.48%  '
